# Geotagger workbook: change the launch-intent button to "Record Location"
#
# - Adds a new "buttonLabel" / "Record Location" column (F) to the survey
#   sheet driving the launch-intent button text.
# - Makes the "survey" sheet the active/selected tab (it was "settings"
#   before), with the selection resting on the newly added F6 cell.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("survey")

# New column F: header + value driving the launch button's label.
$ws1.Range("F1").Value = "buttonLabel"
$ws1.Range("F5").Value = "Record Location"

# Give column F the same width it has in the authored workbook.
$ws1.Columns.Item(6).ColumnWidth = 13.667

# Switch the active sheet/tab from "settings" to "survey" and rest the
# selection on F6, just past the newly-populated column.
$ws1.Activate() | Out-Null
$ws1.Range("F6").Select() | Out-Null
